# "update storage in DG"
# The "Storage" class-diagram slide (SlideID 264) is removed from the deck.
$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 264) {
        $slide.Delete()
    }
}
